# Update "paises.xlsx" (sheet "Pais") with refreshed COVID case counts and
# the new "last updated" timestamp, reproducing the author's commit:
# "Update countries & provincias Spain"
#
# The underlying sheet is kept sorted by "Casos totales" (column B)
# descending. Because several countries' totals changed, some of them
# swapped places with their neighbours. Each affected row below is
# written with its final (post-edit) values, including the correct
# country label for rows whose occupant changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Septiembre de 2020 a las 11:16"

# Filipinas (row 25) - values updated in place
$ws.Range("B25").Value = 232072
$ws.Range("C25").Value = 3714
$ws.Range("D25").Value = 160549
$ws.Range("E25").Value = 67786
$ws.Range("G25").Value = 49
$ws.Range("H25").Value = 3737

# Indonesia (row 26) - values updated in place
$ws.Range("B26").Value = 187537
$ws.Range("C26").Value = 3269
$ws.Range("D26").Value = 134181
$ws.Range("E26").Value = 45524
$ws.Range("G26").Value = 82
$ws.Range("H26").Value = 7832

# Polonia moves ahead of Japon (rows 47/48 swap)
$ws.Range("A47").Value = "Polonia"
$ws.Range("B47").Value = 69820
$ws.Range("C47").Value = 691
$ws.Range("D47").Value = 49820
$ws.Range("E47").Value = 17900
$ws.Range("G47").Value = 8
$ws.Range("H47").Value = 2100

$ws.Range("A48").Value = "Japon"
$ws.Range("B48").Value = 69599
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 59524
$ws.Range("E48").Value = 8756
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 1319

# Singapur (row 52) - values updated in place
$ws.Range("B52").Value = 56948
$ws.Range("C52").Value = 40
$ws.Range("E52").Value = 893

# Croacia moves ahead of Noruega (rows 90/91 swap)
$ws.Range("A90").Value = "Croacia"
$ws.Range("B90").Value = 11428
$ws.Range("C90").Value = 334
$ws.Range("D90").Value = 8530
$ws.Range("E90").Value = 2703
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 195

$ws.Range("A91").Value = "Noruega"
$ws.Range("B91").Value = 11120
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 9348
$ws.Range("E91").Value = 1508
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 264

# Malasia (row 96) - values updated in place
$ws.Range("B96").Value = 9385
$ws.Range("C96").Value = 11
$ws.Range("D96").Value = 9092
$ws.Range("E96").Value = 165

# Hong Kong (row 112) - values updated in place
$ws.Range("B112").Value = 4851
$ws.Range("C112").Value = 12
$ws.Range("D112").Value = 4456
$ws.Range("E112").Value = 301

# Eslovaquia jumps ahead of Ruanda/Surinam/Cuba/Mozambique (rows 118-122 cascade)
$ws.Range("A118").Value = "Eslovaquia"
$ws.Range("B118").Value = 4300
$ws.Range("C118").Value = 137
$ws.Range("D118").Value = 2693
$ws.Range("E118").Value = 1570
$ws.Range("H118").Value = 37

$ws.Range("A119").Value = "Ruanda"
$ws.Range("B119").Value = 4255
$ws.Range("D119").Value = 2163
$ws.Range("E119").Value = 2074
$ws.Range("H119").Value = 18

$ws.Range("A120").Value = "Surinam"
$ws.Range("B120").Value = 4215
$ws.Range("D120").Value = 3318
$ws.Range("E120").Value = 824
$ws.Range("H120").Value = 73

$ws.Range("A121").Value = "Cuba"
$ws.Range("B121").Value = 4214
$ws.Range("D121").Value = 3474
$ws.Range("E121").Value = 640
$ws.Range("H121").Value = 100

$ws.Range("A122").Value = "Mozambique"
$ws.Range("B122").Value = 4207
$ws.Range("D122").Value = 2370
$ws.Range("E122").Value = 1811
$ws.Range("H122").Value = 26

# Sri Lanka (row 128) - values updated in place
$ws.Range("D128").Value = 2907
$ws.Range("E128").Value = 192

# Lituania (row 131) - values updated in place
$ws.Range("B131").Value = 3004
$ws.Range("C131").Value = 26
$ws.Range("D131").Value = 1920
$ws.Range("E131").Value = 998

# Letonia (row 155) - values updated in place
$ws.Range("B155").Value = 1416
$ws.Range("C155").Value = 6
$ws.Range("E155").Value = 194
$ws.Range("G155").Value = 1
$ws.Range("H155").Value = 35

# San Pedro y Miquelon (row 218) - values updated in place
$ws.Range("D218").Value = 5
$ws.Range("E218").Value = 0
